# dijkstra.xlsx: split each of the two input sheets ("dijkstra_input_ccm",
# "dijkstra_input_ram") into a FLASH-only and a CCM-only sheet (columns
# 24/48/72), and append a new "energy" row to each.
#
# xlPasteFormats
$xlPasteFormats = -4122

function Set-TextCell($range, [string]$value, $styleTemplate) {
    # Force the cell to store a genuine text value (not an auto-inferred
    # number) without permanently altering its style: stamp NumberFormat
    # "@" before the write, then paste the original cell's format back on
    # top (keeps the shared style index, e.g. s="1"/s="2").
    $range.NumberFormat = "@"
    $range.Value = $value
    $styleTemplate.Copy()
    $range.PasteSpecial($xlPasteFormats)
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Repurpose the two existing sheets in place (this keeps their
#    sheetId / rId / physical part stable) and insert two brand-new
#    sheets at the correct spots so the final sheetId order comes out
#    1,2,3,4 left-to-right, matching tab order.
# ---------------------------------------------------------------------

# dijkstra_input_ram (sheetId 2) becomes the CCM split of
# dijkstra_input_ccm, and is moved right after dijkstra_input_ccm.
$wb.Worksheets.Item("dijkstra_input_ram").Name = "dijkstra_input_ccm code_CCM"
$moveTarget = $wb.Worksheets.Item("dijkstra_input_ccm")
$wb.Worksheets.Item("dijkstra_input_ccm code_CCM").Move([System.Reflection.Missing]::Value, $moveTarget)

# Two new sheets for the ram FLASH/CCM splits, appended at the end (so
# they're created 3rd and 4th => sheetId 3 and 4).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ramFlash = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ramFlash.Name = "dijkstra_input_ram code_FLASH"
$ramCcm = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item("dijkstra_input_ram code_FLASH"))
$ramCcm.Name = "dijkstra_input_ram code_CCM"

# Finally rename the original ccm sheet to the FLASH split name.
$wb.Worksheets.Item("dijkstra_input_ccm").Name = "dijkstra_input_ccm code_FLASH"

# ---------------------------------------------------------------------
# 2. dijkstra_input_ccm code_FLASH (sheet 1, style s="1")
#    columns 24/48/72 <- old B/D/F (FLASH 24/48/72)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("dijkstra_input_ccm code_FLASH")
$style1 = $ws1.Range("A2")

Set-TextCell $ws1.Range("B1") "24" $style1
Set-TextCell $ws1.Range("C1") "48" $style1
Set-TextCell $ws1.Range("D1") "72" $style1

$ws1.Range("B2").Value = 13658.29965377762
$ws1.Range("C2").Value = 23685.54773270165
$ws1.Range("D2").Value = 31501.07700564517

$ws1.Range("B3").Value = 0.94462
$ws1.Range("C3").Value = 0.5515
$ws1.Range("D3").Value = 0.43392

$ws1.Range("B4").Value = "(742.26, 1686.88)"
$ws1.Range("C4").Value = "(3806.92, 4358.42)"
$ws1.Range("D4").Value = "(6006.38, 6440.3)"

$ws1.Range("E1:G4").Delete()

$ws1.Range("A2").Copy($ws1.Range("A5"))
$ws1.Range("A5").Value = "energy"
$ws1.Range("B5").Value = 42576.27996253967
$ws1.Range("C5").Value = 43106.51259613037
$ws1.Range("D5").Value = 45107.52620315551

# ---------------------------------------------------------------------
# 3. dijkstra_input_ccm code_CCM (sheet 2, style s="2")
#    columns 24/48/72 <- old C/E/G (CCM 24/48/72) of dijkstra_input_ccm
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("dijkstra_input_ccm code_CCM")
$style2 = $ws2.Range("A2")

Set-TextCell $ws2.Range("B1") "24" $style2
Set-TextCell $ws2.Range("C1") "48" $style2
Set-TextCell $ws2.Range("D1") "72" $style2

$ws2.Range("A2").Value = "intensity"
$ws2.Range("B2").Value = 11851.7355026581
$ws2.Range("C2").Value = 22982.01089362742
$ws2.Range("D2").Value = 34024.83088938572

$ws2.Range("A3").Value = "runtime"
$ws2.Range("B3").Value = 0.9448
$ws2.Range("C3").Value = 0.47266
$ws2.Range("D3").Value = 0.31526

$ws2.Range("A4").Value = "timestamp"
$ws2.Range("B4").Value = "(1882.86, 2827.66)"
$ws2.Range("C4").Value = "(4554.48, 5027.14)"
$ws2.Range("D4").Value = "(6636.24, 6951.5)"

$ws2.Range("E1:G4").Delete()

$ws2.Range("A2").Copy($ws2.Range("A5"))
$ws2.Range("A5").Value = "energy"
$ws2.Range("B5").Value = 36951.81501960754
$ws2.Range("C5").Value = 35846.83498764038
$ws2.Range("D5").Value = 35398.00501441955

# ---------------------------------------------------------------------
# 4. dijkstra_input_ram code_FLASH (sheet 3, style s="2")
#    columns 24/48/72 <- old B/D/F (FLASH 24/48/72) of dijkstra_input_ram
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("dijkstra_input_ram code_FLASH")

Set-TextCell $ws3.Range("B1") "24" $style2
Set-TextCell $ws3.Range("C1") "48" $style2
Set-TextCell $ws3.Range("D1") "72" $style2

$style2.Copy($ws3.Range("A2"))
$ws3.Range("A2").Value = "intensity"
$ws3.Range("B2").Value = 14037.15386983452
$ws3.Range("C2").Value = 24262.78369182243
$ws3.Range("D2").Value = 31890.84216007496

$style2.Copy($ws3.Range("A3"))
$ws3.Range("A3").Value = "runtime"
$ws3.Range("B3").Value = 0.9083599999999999
$ws3.Range("C3").Value = 0.53382
$ws3.Range("D3").Value = 0.42204

$style2.Copy($ws3.Range("A4"))
$ws3.Range("A4").Value = "timestamp"
$ws3.Range("B4").Value = "(1110.6, 2018.96)"
$ws3.Range("C4").Value = "(4103.12, 4636.94)"
$ws3.Range("D4").Value = "(6266.48, 6688.52)"

$style2.Copy($ws3.Range("A5"))
$ws3.Range("A5").Value = "energy"
$ws3.Range("B5").Value = 42077.6039943695
$ws3.Range("C5").Value = 42741.46532821655
$ws3.Range("D5").Value = 44415.39638328552

# ---------------------------------------------------------------------
# 5. dijkstra_input_ram code_CCM (sheet 4, style s="2")
#    columns 24/48/72 <- old C/E/G (CCM 24/48/72) of dijkstra_input_ram
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("dijkstra_input_ram code_CCM")

Set-TextCell $ws4.Range("B1") "24" $style2
Set-TextCell $ws4.Range("C1") "48" $style2
Set-TextCell $ws4.Range("D1") "72" $style2

$style2.Copy($ws4.Range("A2"))
$ws4.Range("A2").Value = "intensity"
$ws4.Range("B2").Value = 12112.74827389945
$ws4.Range("C2").Value = 23487.10307937998
$ws4.Range("D2").Value = 34845.6205039571

$style2.Copy($ws4.Range("A3"))
$ws4.Range("A3").Value = "runtime"
$ws4.Range("B3").Value = 0.9086
$ws4.Range("C3").Value = 0.45402
$ws4.Range("D3").Value = 0.30288

$style2.Copy($ws4.Range("A4"))
$ws4.Range("A4").Value = "timestamp"
$ws4.Range("B4").Value = "(2215.02, 3123.62)"
$ws4.Range("C4").Value = "(4832.96, 5286.98)"
$ws4.Range("D4").Value = "(6884.52, 7187.4)"

$style2.Copy($ws4.Range("A5"))
$ws4.Range("A5").Value = "energy"
$ws4.Range("B5").Value = 36318.62216949463
$ws4.Range("C5").Value = 35189.92798233031
$ws4.Range("D5").Value = 34828.33707618713

# ---------------------------------------------------------------------
# 6. Sheet view niceties matching the original workbook (first sheet
#    tab-selected/active).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("dijkstra_input_ccm code_FLASH").Activate()
